# Fruta / hortaliza, semanal
# Insert this week's new price record at the top of the data table
# (row 221, just below the header + row-220 block), pushing every
# existing record down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 221:235 down to 222:236, duplicating row 221's formatting
# (e.g. the date style on column D) onto the freshly inserted row.
$ws.Rows.Item(221).Insert()

# Populate the new weekly record.
$ws.Range("A221").Value = 10
$ws.Range("B221").Value = "Vega Modelo de Temuco"
$ws.Range("C221").Value = "La Araucanía"
$ws.Range("D221").Value = 44714
$ws.Range("E221").Value = 9
$ws.Range("F221").Value = "Fruta"
$ws.Range("G221").Value = 100102
$ws.Range("H221").Value = "Cítricos"
$ws.Range("I221").Value = 100102006
$ws.Range("J221").Value = "Pomelo"
$ws.Range("K221").Value = "Start Ruby"
$ws.Range("L221").Value = "Primera"
$ws.Range("M221").Value = 90
$ws.Range("N221").Value = 12000
$ws.Range("O221").Value = 12000
$ws.Range("P221").Value = 12000
$ws.Range("Q221").Value = "$/bandeja 15 kilos granel"
$ws.Range("R221").Value = "Región de O'Higgins"
$ws.Range("S221").Value = 800
$ws.Range("T221").Value = 15
